$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 1.524029
$ws.Range("H2").Value = 3.048058
$ws.Range("I2").Value = 0.09030204154573296
$ws.Range("J2").Value = 0.06866669168778029
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.1750253333333333
$ws.Range("N2").Value = 0.525076
$ws.Range("O2").Value = 0.655128508251241
$ws.Range("P2").Value = 0.7402223729081171
$ws.Range("Q2").Value = 0.2667436837346667
$ws.Range("R2").Value = 1.600462102408
$ws.Range("S2").Value = 0.05915944176989762
$ws.Range("T2").Value = 0.05082862146087881

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 1.524029
$ws.Range("H3").Value = 3.048058
$ws.Range("I3").Value = 0.09030204154573296
$ws.Range("J3").Value = 0.06866669168778029
$ws.Range("K3").Value = 2
$ws.Range("M3").Value = 0.09213650000000001
$ws.Range("N3").Value = 0.184273
$ws.Range("O3").Value = 0.3448714917487591
$ws.Range("P3").Value = 0.2597776270918828
$ws.Range("Q3").Value = 0.1404186979585
$ws.Range("R3").Value = 0.561674791834
$ws.Range("S3").Value = 0.03114259977583534
$ws.Range("T3").Value = 0.01783807022690148

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.502875000000001
$ws.Range("H4").Value = 25.508625
$ws.Range("I4").Value = 0.5038138851085998
$ws.Range("J4").Value = 0.5746586476550659
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1750253333333333
$ws.Range("N4").Value = 0.525076
$ws.Range("O4").Value = 0.655128508251241
$ws.Range("P4").Value = 0.7402223729081171
$ws.Range("Q4").Value = 1.488218531166667
$ws.Range("R4").Value = 13.3939667805
$ws.Range("S4").Value = 0.3300628389874591
$ws.Range("T4").Value = 0.4253751877794025

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 8.502875000000001
$ws.Range("H5").Value = 25.508625
$ws.Range("I5").Value = 0.5038138851085998
$ws.Range("J5").Value = 0.5746586476550659
$ws.Range("K5").Value = 2
$ws.Range("M5").Value = 0.09213650000000001
$ws.Range("N5").Value = 0.184273
$ws.Range("O5").Value = 0.3448714917487591
$ws.Range("P5").Value = 0.2597776270918828
$ws.Range("Q5").Value = 0.7834251424375002
$ws.Range("R5").Value = 4.700550854625001
$ws.Range("S5").Value = 0.1737510461211407
$ws.Range("T5").Value = 0.1492834598756634

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.3890603333333333
$ws.Range("H6").Value = 1.167181
$ws.Range("I6").Value = 0.02305267313447669
$ws.Range("J6").Value = 0.02629426929239375
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.1750253333333333
$ws.Range("N6").Value = 0.525076
$ws.Range("O6").Value = 0.655128508251241
$ws.Range("P6").Value = 0.7402223729081171
$ws.Range("Q6").Value = 0.06809541452844445
$ws.Range("R6").Value = 0.6128587307559999
$ws.Range("S6").Value = 0.01510246336179318
$ws.Range("T6").Value = 0.01946360640950074

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.3890603333333333
$ws.Range("H7").Value = 1.167181
$ws.Range("I7").Value = 0.02305267313447669
$ws.Range("J7").Value = 0.02629426929239375
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 0.09213650000000001
$ws.Range("N7").Value = 0.184273
$ws.Range("O7").Value = 0.3448714917487591
$ws.Range("P7").Value = 0.2597776270918828
$ws.Range("Q7").Value = 0.03584665740216667
$ws.Range("R7").Value = 0.215079944413
$ws.Range("S7").Value = 0.007950209772683519
$ws.Range("T7").Value = 0.006830662882893009

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.4482056666666667
$ws.Range("H8").Value = 1.344617
$ws.Range("I8").Value = 0.02655716310671665
$ws.Range("J8").Value = 0.0302915498908315
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.1750253333333333
$ws.Range("N8").Value = 0.525076
$ws.Range("O8").Value = 0.655128508251241
$ws.Range("P8").Value = 0.7402223729081171
$ws.Range("Q8").Value = 0.07844734621022223
$ws.Range("R8").Value = 0.7060261158919999
$ws.Range("S8").Value = 0.01739835464948817
$ws.Range("T8").Value = 0.02242248293925591

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.4482056666666667
$ws.Range("H9").Value = 1.344617
$ws.Range("I9").Value = 0.02655716310671665
$ws.Range("J9").Value = 0.0302915498908315
$ws.Range("K9").Value = 2
$ws.Range("M9").Value = 0.09213650000000001
$ws.Range("N9").Value = 0.184273
$ws.Range("O9").Value = 0.3448714917487591
$ws.Range("P9").Value = 0.2597776270918828
$ws.Range("Q9").Value = 0.04129610140683334
$ws.Range("R9").Value = 0.247776608441
$ws.Range("S9").Value = 0.00915880845722848
$ws.Range("T9").Value = 0.007869066951575588

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.295005
$ws.Range("H10").Value = 3.885015
$ws.Range("I10").Value = 0.07673187013628475
$ws.Range("J10").Value = 0.08752167025935917
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.1750253333333333
$ws.Range("N10").Value = 0.525076
$ws.Range("O10").Value = 0.655128508251241
$ws.Range("P10").Value = 0.7402223729081171
$ws.Range("Q10").Value = 0.2266586817933333
$ws.Range("R10").Value = 2.03992813614
$ws.Range("S10").Value = 0.05026923561771218
$ws.Range("T10").Value = 0.06478549844026463

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.295005
$ws.Range("H11").Value = 3.885015
$ws.Range("I11").Value = 0.07673187013628475
$ws.Range("J11").Value = 0.08752167025935917
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 0.09213650000000001
$ws.Range("N11").Value = 0.184273
$ws.Range("O11").Value = 0.3448714917487591
$ws.Range("P11").Value = 0.2597776270918828
$ws.Range("Q11").Value = 0.1193172281825
$ws.Range("R11").Value = 0.7159033690950001
$ws.Range("S11").Value = 0.02646263451857258
$ws.Range("T11").Value = 0.02273617181909454

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 4.717841
$ws.Range("H12").Value = 9.435682
$ws.Range("I12").Value = 0.2795423669681891
$ws.Range("J12").Value = 0.2125671712145694
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.1750253333333333
$ws.Range("N12").Value = 0.525076
$ws.Range("O12").Value = 0.655128508251241
$ws.Range("P12").Value = 0.7402223729081171
$ws.Range("Q12").Value = 0.8257416936386667
$ws.Range("R12").Value = 4.954450161832
$ws.Range("S12").Value = 0.1831361738648907
$ws.Range("T12").Value = 0.1573469758788146

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 4.717841
$ws.Range("H13").Value = 9.435682
$ws.Range("I13").Value = 0.2795423669681891
$ws.Range("J13").Value = 0.2125671712145694
$ws.Range("K13").Value = 2
$ws.Range("M13").Value = 0.09213650000000001
$ws.Range("N13").Value = 0.184273
$ws.Range("O13").Value = 0.3448714917487591
$ws.Range("P13").Value = 0.2597776270918828
$ws.Range("Q13").Value = 0.4346853572965
$ws.Range("R13").Value = 1.738741429186
$ws.Range("S13").Value = 0.05522019533575483
$ws.Range("T13").Value = 0.05522019533575483
